# release juni 2022: nieuwe wijkindelingen (+correcties)
# Update v9900_type_ggw7 (col D) from 2.00 -> 1.00 and, where applicable,
# v9900_eerstewijk (col E) text values for a set of "gemeente" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row -> new E value (column 5), or $null when E is unchanged
$changes = @{
    17  = "110291"
    41  = $null
    94  = $null
    97  = $null
    107 = $null
    116 = $null
    121 = $null
    138 = "24045"
    162 = "31040AARBU"
    181 = "34013AW"
    200 = "36010LEDBU"
    207 = "37010GIN"
    211 = "37017INDU"
    238 = "42028AD"
    240 = "43005BAL"
    289 = $null
    295 = "71070BER"
}

foreach ($row in $changes.Keys) {
    # column D: v9900_type_ggw7, 2.00 -> 1.00
    $ws.Cells.Item($row, 4).Value = 1.00

    $newE = $changes[$row]
    if ($newE -ne $null) {
        $cell = $ws.Cells.Item($row, 5)
        if ($newE -match '^[0-9]+$') {
            # Purely-numeric replacement text (e.g. "110291") must stay a
            # text value, not be auto-coerced to a number. Prefix with an
            # apostrophe to force text entry, then reset the cell style so
            # no stray "Text" number-format style gets attached.
            $cell.Value = "'" + $newE
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newE
        }
    }
}
